$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Metadata sheet -------------------------------------------------------
# A brand new row ("Jurisdiction" / "iso:code:3166:FR") is inserted as the
# new row 11; every existing row from the old row 11 ("Description") down to
# the old row 21 ("Count") shifts down by one (to rows 12..22).
#
# We shift the values ourselves (row by row, bottom-up so a source row is
# never clobbered before it is read) rather than using Range.Insert, because
# every row 2..21 already carries the shared "s=2" body style - plain .Value
# assignment leaves that formatting untouched and does not add any new style
# entries to styles.xml. Only the brand-new row 22 needs its format fixed up
# afterwards (it has no prior cell/style to inherit).

for ($r = 21; $r -ge 11; $r--) {
    $destRow = $r + 1
    $aVal = $ws.Cells.Item($r, 1).Value()
    $bVal = $ws.Cells.Item($r, 2).Value()

    $ws.Cells.Item($destRow, 1).Value = $aVal
    if ($bVal -eq $null) {
        $ws.Cells.Item($destRow, 2).Value = ""
    } else {
        $ws.Cells.Item($destRow, 2).Value = $bVal
    }
}

# New row 11.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"

# Other metadata text updates (Version, Date) for the 0.2.0 publication.
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# The "Case Sensitive" value ("true", now on row 15) and the "Count" value
# ("1", now on row 22) must stay plain text - matching the original
# shared-string cells - rather than being auto-typed into a Boolean/Number by
# a bare .Value assignment. Use a leading apostrophe to force text, then wipe
# the resulting "quote prefix" cell format by re-pasting the formatting from
# an untouched neighbouring body cell, so the final style matches the rest of
# the table exactly.
$ws.Range("B15").Value = "'true"
$ws.Range("B16").Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B22").Value = "'1"
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
